$d = $word.ActiveDocument

# Update the date/title line at the top of the document
$d.Content.Find.Execute("2025-12-04 Thursday", $false, $false, $false, $false, $false, $true, 1, $false, "2025-12-05 Friday", 2) | Out-Null

# Update the division problems in the table. Each replacement is scoped
# to its specific cell and uses wdReplaceOne (1) rather than wdReplaceAll (2)
# because some expressions (e.g. "44÷6=") occur more than once in the
# document, and ReplaceAll would match/replace every occurrence even when
# the Find is confined to a single cell Range.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Find.Execute("77÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "37÷6=", 1) | Out-Null
$t.Cell(1, 2).Range.Find.Execute("94÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "56÷6=", 1) | Out-Null
$t.Cell(1, 3).Range.Find.Execute("44÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "12÷7=", 1) | Out-Null
$t.Cell(1, 4).Range.Find.Execute("19÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "89÷8=", 1) | Out-Null
$t.Cell(1, 5).Range.Find.Execute("33÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "99÷8=", 1) | Out-Null
$t.Cell(5, 1).Range.Find.Execute("54÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "81÷2=", 1) | Out-Null
$t.Cell(5, 2).Range.Find.Execute("53÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "26÷4=", 1) | Out-Null
$t.Cell(5, 3).Range.Find.Execute("77÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "56÷3=", 1) | Out-Null
$t.Cell(5, 4).Range.Find.Execute("91÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "27÷2=", 1) | Out-Null
$t.Cell(5, 5).Range.Find.Execute("16÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "96÷3=", 1) | Out-Null
$t.Cell(9, 1).Range.Find.Execute("44÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "69÷4=", 1) | Out-Null
$t.Cell(9, 2).Range.Find.Execute("35÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "89÷7=", 1) | Out-Null
$t.Cell(9, 3).Range.Find.Execute("19÷2=", $false, $false, $false, $false, $false, $true, 1, $false, "12÷6=", 1) | Out-Null
$t.Cell(9, 4).Range.Find.Execute("43÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "44÷8=", 1) | Out-Null
$t.Cell(9, 5).Range.Find.Execute("86÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "17÷7=", 1) | Out-Null
$t.Cell(13, 1).Range.Find.Execute("93÷5=", $false, $false, $false, $false, $false, $true, 1, $false, "74÷4=", 1) | Out-Null
$t.Cell(13, 2).Range.Find.Execute("52÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "59÷4=", 1) | Out-Null
$t.Cell(13, 3).Range.Find.Execute("77÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "93÷2=", 1) | Out-Null
$t.Cell(13, 4).Range.Find.Execute("22÷2=", $false, $false, $false, $false, $false, $true, 1, $false, "17÷5=", 1) | Out-Null
$t.Cell(13, 5).Range.Find.Execute("55÷5=", $false, $false, $false, $false, $false, $true, 1, $false, "40÷4=", 1) | Out-Null
$t.Cell(17, 1).Range.Find.Execute("50÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "46÷6=", 1) | Out-Null
$t.Cell(17, 2).Range.Find.Execute("71÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "68÷5=", 1) | Out-Null
$t.Cell(17, 3).Range.Find.Execute("62÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "70÷5=", 1) | Out-Null
$t.Cell(17, 4).Range.Find.Execute("27÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "54÷8=", 1) | Out-Null
$t.Cell(17, 5).Range.Find.Execute("76÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "63÷3=", 1) | Out-Null
